$d = $word.ActiveDocument
$enDash = [char]0x2013

# ---------------------------------------------------------------------------
# 1. Merge the three runs "Felipe Cardoso " + "-" + " Groceries Store App"
#    into a single run "Felipe Cardoso - Groceries Store App".
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.ClearFormatting()
$felipeText = "Felipe Cardoso " + $enDash + " Groceries Store App"
$findRng.Find.Execute($felipeText, $true, $false, $false, $false, $false, $true, 1, $false, $felipeText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. The very last paragraph in the body is an empty "ListParagraph" item
#    (pStyle ListParagraph + numPr). Strip its numbering/style so it becomes
#    a bare empty paragraph, then append the new "Camron Darpoh" block of
#    paragraphs after it.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$lastRange.ListFormat.RemoveNumbers()
$lastPara.Style = "Normal"

# Create a fresh trailing paragraph to anchor the XML insertion point right
# after the (now bare) last paragraph.
$lastRange.InsertParagraphAfter()
$anchorPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertionPoint = $anchorPara.Range
$insertionPoint.Collapse(1) | Out-Null

$newContentXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:t>Camron Darpoh &#8211; Events app</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>App that would display all university events and programmes</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Display all community and student events</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Much easier to spread information as opposed to leaflets and </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:t>brunel</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:t xml:space="preserve"> website</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:left="360"/>
            </w:pPr>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>Reason not picked:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t>Not sure how applicable it would be whether students would be willing to download</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Very demanding as a lot of communities run separately </w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($newContentXml)

# Remove the helper trailing empty paragraph that was created purely to
# anchor the insertion point (the inserted XML's last paragraph is now the
# true last paragraph of the body).
$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$cleanupRange = $d.Range($trailingPara.Range.Start - 1, $trailingPara.Range.End)
$cleanupRange.Delete()
